$wb = $excel.ActiveWorkbook

# Update status from "Ready for handoff" to "In Translation" for the
# 41797f57-...md and 4858cf89-...md files (rows 8 and 9) across the
# Overview, zh-cn and de-de sheets, as part of regenerating the report.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("C9").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("C9").Value = "In Translation"
